$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(6)

# Work right-to-left over the original run layout so earlier (lower) character
# offsets stay valid while later ones are edited first.

# Run 6 (original): "), ODE models"  -> "), ODE models "
$r6 = $para.Characters(225, 13)
$r6.Text = "), ODE models "

# Run 4 (original): "speciality" -> "Ms"  (keeps its own rPr, incl. err="1")
$r4 = $para.Characters(214, 10)
$r4.Text = "Ms"

# Insert a brand-new run right after the (now 2-char) run 4, carrying the
# text that used to continue inside runs 3/4 plus new content.
$r4 = $para.Characters(214, 2)
$r4.InsertAfter(" x1, machine learning (not my specialty")

# Run 3 (original): " models, autoregressive models (time series data), shiny apps, power analyses, non-linear models, machine learning (not my "
# -> " models, autoregressive models (time series data x2), shiny apps, power analyses x2, non-linear models & GAM/"
$r3 = $para.Characters(91, 123)
$r3.Text = " models, autoregressive models (time series data x2), shiny apps, power analyses x2, non-linear models & GAM/"

# Run 1 (original): "Possible options: zero-inflated models and bias-reduction models, correlated data e.g. "
# -> "Possible options: zero-inflated models/hurdle models (x1) and bias-reduction models, correlated data e.g. "
$r1 = $para.Characters(1, 87)
$r1.Text = "Possible options: zero-inflated models/hurdle models (x1) and bias-reduction models, correlated data e.g. "
